$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 86, pushing the existing rows 86-159 down to 87-160.
$ws.Rows("86:86").Insert()

# Populate the newly inserted row 86 with the new record.
$ws.Cells.Item(86, 1).Value = 5
$ws.Cells.Item(86, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(86, 3).Value = "Maule"
$ws.Cells.Item(86, 4).Value = 44566
$ws.Cells.Item(86, 5).Value = 7
$ws.Cells.Item(86, 6).Value = 100112021
$ws.Cells.Item(86, 7).Value = "Ají"
$ws.Cells.Item(86, 8).Value = "Americana (o)"
$ws.Cells.Item(86, 9).Value = "Primera"
$ws.Cells.Item(86, 10).Value = 150
$ws.Cells.Item(86, 11).Value = 17000
$ws.Cells.Item(86, 12).Value = 17000
$ws.Cells.Item(86, 13).Value = 17000
$ws.Cells.Item(86, 14).Value = "$/caja 14 kilos"
$ws.Cells.Item(86, 15).Value = "Región del Maule"
$ws.Cells.Item(86, 16).Value = 1214
$ws.Cells.Item(86, 17).Value = 14
$ws.Cells.Item(86, 18).Value = "Hortaliza"

# Make sure column D keeps the date-formatted style used throughout the table.
$ws.Cells.Item(86, 4).NumberFormat = $ws.Cells.Item(87, 4).NumberFormat
